# "update to tasks list"
# The "Integration" sheet's Assigned (F) / Status (G) placeholders ("??")
# are filled in with the actual assignees, the "Assigned" column is widened
# back down, the sheet is frozen at C2 with a 120% zoom, and a couple of
# column widths are tweaked.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Integration")

# --- Fill in the placeholder "??" assignments with real names/status ---
$ws.Range("F8").Value2  = "Sushma"
$ws.Range("G8").Value2  = "In progress"
$ws.Range("F10").Value2 = "Sushma"
$ws.Range("F13").Value2 = "Henry"
$ws.Range("F14").Value2 = "Henry"
$ws.Range("F26").Value2 = "Scott, Sushma"
$ws.Range("F27").Value2 = "Scott"
$ws.Range("F28").Value2 = "Team"

# --- Column width tweaks (Detail / Status columns) ---
$ws.Columns.Item(3).ColumnWidth = 49.33
$ws.Columns.Item(7).ColumnWidth = 13.33

# --- Freeze panes at C2 (keep Category/Step visible), zoom to 120% ---
$ws.Activate()
$ws.Range("C2").Select()
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.Zoom = 120
